$wb = $excel.ActiveWorkbook

# "Ready for handoff" lives in the shared-strings table and is referenced by
# several cells (Overview!E2:F3, zh-cn!C2:C3, de-de!C2:C3). Updating the text
# on each matching cell updates the shared string itself, so every reference
# repaints to the new text.
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    $cols = $used.Columns.Count
    $touchedCols = @{}
    for ($r = 1; $r -le $rows; $r++) {
        for ($c = 1; $c -le $cols; $c++) {
            $cell = $used.Cells.Item($r, $c)
            # NOTE: literal goes on the left of -eq so PowerShell's type
            # coercion follows the string's type (not e.g. a boolean cell's
            # type, which would otherwise coerce "Ready for handoff" to
            # $true and false-match every True/False cell).
            if ("Ready for handoff" -eq $cell.Value()) {
                $cell.Value = "In Translation"
                $touchedCols[$cell.Column] = $true
            }
        }
    }
    # Shrink only the columns whose text actually changed (from the 17-char
    # "Ready for handoff" down to the 14-char "In Translation"), matching a
    # regenerated report's column-width refresh. Leave every other column's
    # width untouched.
    foreach ($colIndex in $touchedCols.Keys) {
        $ws.Cells.Item(1, $colIndex).EntireColumn.ColumnWidth = 12.5
    }
}

